## all figures and tables updated
# The "table 1" sheet (sheet4) is the active sheet in this workbook.
# Its data rows (5-11) need to be reversed in order (ENS..COL top-to-bottom
# instead of COL..ENS), and the current selection moves from D17 to G7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E")
$pairs = @(@(5, 11), @(6, 10), @(7, 9))

foreach ($pair in $pairs) {
    $top = $pair[0]
    $bottom = $pair[1]
    foreach ($col in $cols) {
        $topCell = $ws.Range($col + $top)
        $bottomCell = $ws.Range($col + $bottom)
        $tmp = $topCell.Value2
        $topCell.Value2 = $bottomCell.Value2
        $bottomCell.Value2 = $tmp
    }
}

[void]$ws.Range("G7").Select()
